# Updated symbol list (prices / 1h volume changes, plus a few coins that
# shifted rank and swapped rows) to match the latest coinranking.com pull.
# Numeric-looking strings are written with a leading apostrophe so they
# stay text (matching column D "Price" / column E "Volume(1h)" which are
# stored as text, not numbers, in this sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'327.07"
$ws.Range("E2").Value = "'-0.75%"
$ws.Range("D3").Value = "'43.82"
$ws.Range("E3").Value = "'0.66%"
$ws.Range("D4").Value = "'5.541"
$ws.Range("E4").Value = "'-0.60%"
$ws.Range("E5").Value = "'-2.17%"
$ws.Range("D6").Value = "'1.897"
$ws.Range("E6").Value = "'0.74%"
$ws.Range("D7").Value = "'4.262"
$ws.Range("E7").Value = "'-2.53%"
$ws.Range("D8").Value = "'0.9446"
$ws.Range("E8").Value = "'0.18%"
$ws.Range("D9").Value = "'2.539"
$ws.Range("E9").Value = "'-9.48%"
$ws.Range("D10").Value = "'0.1176"
$ws.Range("E10").Value = "'-1.03%"
$ws.Range("D11").Value = "'0.1836"
$ws.Range("E11").Value = "'-4.09%"
$ws.Range("B12").Value = "MCDex"
$ws.Range("C12").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D12").Value = "'10.08"
$ws.Range("E12").Value = "'15.47%"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "'0.09642"
$ws.Range("E13").Value = "'-0.54%"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.04447"
$ws.Range("E14").Value = "'2.91%"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.1066"
$ws.Range("E15").Value = "'-0.34%"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001290"
$ws.Range("E16").Value = "'-0.27%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.005937"
$ws.Range("E17").Value = "'-0.47%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.407"
$ws.Range("E18").Value = "'-3.55%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3444"
$ws.Range("E19").Value = "'-2.61%"
$ws.Range("D20").Value = "'0.1406"
$ws.Range("E20").Value = "'2.67%"
$ws.Range("D21").Value = "'0.2510"
$ws.Range("E21").Value = "'0.60%"
$ws.Range("D22").Value = "'0.04191"
$ws.Range("E22").Value = "'-4.46%"
$ws.Range("D23").Value = "'0.001249"
$ws.Range("E23").Value = "'0.76%"
$ws.Range("D24").Value = "'0.004285"
$ws.Range("E24").Value = "'-1.35%"
$ws.Range("E25").Value = "'2.26%"
$ws.Range("D26").Value = "'0.0003999"
$ws.Range("E26").Value = "'-0.20%"
$ws.Range("D38").Value = "'0.02630"
$ws.Range("E38").Value = "'-5.03%"
$ws.Range("D39").Value = "'0.05475"
$ws.Range("E39").Value = "'-3.91%"
$ws.Range("D40").Value = "'0.007581"
$ws.Range("E40").Value = "'-4.84%"
$ws.Range("E41").Value = "'-2.05%"
$ws.Range("D42").Value = "'0.008183"
$ws.Range("E42").Value = "'-16.07%"
$ws.Range("D43").Value = "'0.002006"
$ws.Range("E43").Value = "'-4.77%"
$ws.Range("D44").Value = "'0.008801"
$ws.Range("E44").Value = "'-12.38%"
$ws.Range("D45").Value = "'0.00007097"
$ws.Range("E45").Value = "'-3.03%"
$ws.Range("E46").Value = "'-0.20%"
$ws.Range("B47").Value = "CoinbaseStockToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D47").Value = "'0.002276"
$ws.Range("E47").Value = "'-0.20%"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "'0.003611"
$ws.Range("E48").Value = "'4.74%"
$ws.Range("D49").Value = "'0.00002105"
$ws.Range("E49").Value = "'-0.20%"
$ws.Range("D50").Value = "'0.0002005"
$ws.Range("E50").Value = "'-0.20%"
